$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 4 sending clusters x 3 target clusters (Resolving-Mac added as a new
# sending cluster; it is not a target cluster in this edit).
$data = @(
    @("ECs", "Vtn", "Itgb8", "ECs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 3.987076, 11.961228, 0.2813308272685638, 0.2813308272685638, 11.034440157952, 99.30996142156799, 0.01324081695213191, 0.01324081695213191),
    @("ECs", "Vtn", "Itgb8", "FAPs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 10.131229, 30.393687, 0.7148664925918803, 0.7148664925918804, 28.03870308140799, 252.348327732672, 0.03364514463459698, 0.03364514463459698),
    @("ECs", "Vtn", "Itgb8", "MuSCs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.00380268013955587, 0.00380268013955587, 0.1491498349013333, 1.342348514112, 0.0001789728916102787, 0.0001789728916102787),
    @("FAPs", "Vtn", "Itgb8", "ECs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 3.987076, 11.961228, 0.2813308272685638, 0.2813308272685638, 83.94060212117334, 755.46541909056, 0.1007248334875628, 0.1007248334875628),
    @("FAPs", "Vtn", "Itgb8", "FAPs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 10.131229, 30.393687, 0.7148664925918803, 0.7148664925918804, 213.2945202166933, 1919.65068195024, 0.2559435420968568, 0.2559435420968568),
    @("FAPs", "Vtn", "Itgb8", "MuSCs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.00380268013955587, 0.00380268013955587, 1.134604635004444, 10.21144171504, 0.001361472994559479, 0.00136147299455948),
    @("MuSCs", "Vtn", "Itgb8", "ECs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 3.987076, 11.961228, 0.2813308272685638, 0.2813308272685638, 139.4576025853627, 1255.118423268264, 0.1673426618826, 0.1673426618826),
    @("MuSCs", "Vtn", "Itgb8", "FAPs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 10.131229, 30.393687, 0.7148664925918803, 0.7148664925918804, 354.3641775535007, 3189.277597981506, 0.4252205949929702, 0.4252205949929703),
    @("MuSCs", "Vtn", "Itgb8", "MuSCs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.00380268013955587, 0.00380268013955587, 1.885014382569556, 16.965129443126, 0.002261929924351674, 0.002261929924351674),
    @("Resolving-Mac", "Vtn", "Itgb8", "ECs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 3.987076, 11.961228, 0.2813308272685638, 0.2813308272685638, 0.018763179656, 0.168868616904, 0.00002251494626902504, 0.00002251494626902504),
    @("Resolving-Mac", "Vtn", "Itgb8", "FAPs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 10.131229, 30.393687, 0.7148664925918803, 0.7148664925918804, 0.047677563674, 0.429098073066, 0.00005721086745629834, 0.00005721086745629836),
    @("Resolving-Mac", "Vtn", "Itgb8", "MuSCs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.00380268013955587, 0.00380268013955587, 0.0002536173206666666, 0.002282555886, 0.000000304329034438367, 0.000000304329034438367)
)

$rowIdx = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $val
        $colIdx++
    }
    $rowIdx++
}
